# Fruta / hortaliza, semanal
# Update the weekly price data: dates (col D) and associated measurement
# columns (I, J, K, L, M, O, P) are refreshed with this week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44204
$ws.Range("J2").Value = 7000

# Row 3
$ws.Range("D3").Value = 44214
$ws.Range("J3").Value = 7000

# Row 4
$ws.Range("D4").Value = 44188
$ws.Range("J4").Value = 12000
$ws.Range("K4").Value = 3000
$ws.Range("M4").Value = 3000
$ws.Range("P4").Value = 30

# Row 5
$ws.Range("D5").Value = 44160
$ws.Range("J5").Value = 7000

# Row 6
$ws.Range("D6").Value = 44232
$ws.Range("J6").Value = 16000

# Row 8
$ws.Range("D8").Value = 44162
$ws.Range("J8").Value = 7000

# Row 9
$ws.Range("D9").Value = 44229
$ws.Range("J9").Value = 16000

# Row 10
$ws.Range("D10").Value = 44189
$ws.Range("J10").Value = 16000

# Row 11
$ws.Range("D11").Value = 44181
$ws.Range("J11").Value = 12000

# Row 12
$ws.Range("D12").Value = 44231
$ws.Range("J12").Value = 12000
$ws.Range("K12").Value = 3000
$ws.Range("M12").Value = 3000
$ws.Range("P12").Value = 30

# Row 13
$ws.Range("D13").Value = 44230
$ws.Range("J13").Value = 16000

# Row 14
$ws.Range("D14").Value = 44159
$ws.Range("J14").Value = 7000
$ws.Range("O14").Value = "Provincia de Chacabuco"

# Row 15
$ws.Range("I15").Value = "Primera"
$ws.Range("J15").Value = 9000
$ws.Range("K15").Value = 3000
$ws.Range("L15").Value = 3000
$ws.Range("M15").Value = 3000
$ws.Range("P15").Value = 30

# Row 16
$ws.Range("D16").Value = 44245
$ws.Range("I16").Value = "Segunda"
$ws.Range("J16").Value = 5000
$ws.Range("K16").Value = 2500
$ws.Range("L16").Value = 2500
$ws.Range("M16").Value = 2500
$ws.Range("O16").Value = "Región Metropolitana"
$ws.Range("P16").Value = 25

# Row 17
$ws.Range("D17").Value = 44187
$ws.Range("J17").Value = 12000

# Row 18
$ws.Range("D18").Value = 44209
$ws.Range("K18").Value = 2500
$ws.Range("M18").Value = 2750
$ws.Range("P18").Value = 28

# Row 19
$ws.Range("D19").Value = 44186
$ws.Range("J19").Value = 10000

# Row 20
$ws.Range("D20").Value = 44168

# Row 21
$ws.Range("D21").Value = 44215
$ws.Range("J21").Value = 16000

# Row 22
$ws.Range("D22").Value = 44210
$ws.Range("J22").Value = 8800
$ws.Range("K22").Value = 2500
$ws.Range("M22").Value = 2750
$ws.Range("P22").Value = 28

# Row 23
$ws.Range("D23").Value = 44161
$ws.Range("J23").Value = 7000

# Row 24
$ws.Range("D24").Value = 44167
